$wb = $excel.ActiveWorkbook

# Sheet "展览" and "全部类型" share identical data rows; update column F (想去人数)
# for the rows that changed.
$changes = @{
    4  = 260
    5  = 43
    11 = 4233
    17 = 59
    18 = 2972
    19 = 57
    20 = 419
    22 = 15
    23 = 65
    28 = 46
    29 = 192
    31 = 392
    32 = 1667
    33 = 241
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $changes.Keys) {
        $ws.Cells.Item($row, 6).Value = $changes[$row]
    }
}
